# Updates cryptos list values (Price and Volume(1h) columns) per the Aug 8 2023 GitHub Actions refresh.
# D-column values are prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the source data's text-formatted numbers, e.g. thousand separators like '29.179.78').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.179.78"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "'1.830.78"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'242.54"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'0.6237"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "'0.07386"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "'0.2936"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "'23.17"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "'0.07675"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'1.830.01"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "'4.965"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "'0.6678"
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").Value = "'82.66"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'0.000009020"
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").Value = "'5.867"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "'29.147.93"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'2.074.43"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("D20").Value = "'235.73"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Value = "'12.46"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "'7.404"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'158.40"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").Value = "'0.1421"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").Value = "'8.537"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "'17.65"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").Value = "'1.486"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'0.05844"
$ws.Range("E30").Value = "  +5.11%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.088"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.098"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "'1.208"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "'1.870"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "'0.7328"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "'1.143"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "'2.605"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'1.226.51"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "'0.01757"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").Value = "'6.293"
$ws.Range("E41").Value = "  -4.77%  "
$ws.Range("D42").Value = "'0.9198"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "'102.09"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'1.976.31"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "'65.23"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5045"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000118"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").Value = "'9.151"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'0.4030"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "'0.1129"
$ws.Range("E51").Value = "  +2.20%  "
